$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.303.23'
$ws.Range("E2").Value = '  +0.31%  '
$ws.Range("D3").Value = '3.499.82'
$ws.Range("E3").Value = '  -0.59%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = "'590.96"
$ws.Range("E5").Value = '  +0.75%  '
$ws.Range("D6").Value = "'133.75"
$ws.Range("E6").Value = '  -0.19%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  -0.46%  '
$ws.Range("E9").Value = '  +3.39%  '
$ws.Range("E10").Value = '  +0.02%  '
$ws.Range("D11").Value = "'0.385"
$ws.Range("E11").Value = '  +2.30%  '
$ws.Range("D12").Value = '4.097.48'
$ws.Range("E14").Value = '  +1.07%  '
$ws.Range("D15").Value = '3.502.14'
$ws.Range("E15").Value = '  -0.50%  '
$ws.Range("D16").Value = '64.381.18'
$ws.Range("E16").Value = '  +0.36%  '
$ws.Range("E17").Value = '  -6.43%  '
$ws.Range("D18").Value = "'9.85"
$ws.Range("E18").Value = '  +0.80%  '
$ws.Range("E19").Value = '  +2.49%  '
$ws.Range("E20").Value = '  -2.55%  '
$ws.Range("D21").Value = "'392.98"
$ws.Range("E21").Value = '  +2.59%  '
$ws.Range("E22").Value = '  +0.68%  '
$ws.Range("D23").Value = '3.640.17'
$ws.Range("D24").Value = "'74.62"
$ws.Range("E24").Value = '  +0.80%  '
$ws.Range("E25").Value = '  -0.15%  '
$ws.Range("E26").Value = '  +0.23%  '
$ws.Range("E27").Value = '  +0.14%  '
$ws.Range("D28").Value = "'7.36"
$ws.Range("E28").Value = '  -1.54%  '
$ws.Range("E29").Value = '  +1.50%  '
$ws.Range("D30").Value = "'8.21"
$ws.Range("E30").Value = '  -2.36%  '
$ws.Range("E31").Value = '  -6.80%  '
$ws.Range("D32").Value = '3.520.16'
$ws.Range("E32").Value = '  -0.40%  '
$ws.Range("D33").Value = "'0.153"
$ws.Range("E34").Value = '  +0.01%  '
$ws.Range("D35").Value = "'23.44"
$ws.Range("E35").Value = '  -0.49%  '
$ws.Range("D36").Value = "'5.15"
$ws.Range("E36").Value = '  -4.67%  '
$ws.Range("D37").Value = "'6.88"
$ws.Range("E37").Value = '  -1.11%  '
$ws.Range("D39").Value = "'167.16"
$ws.Range("E39").Value = '  +5.26%  '
$ws.Range("E40").Value = '  -0.93%  '
$ws.Range("D41").Value = "'0.808"
$ws.Range("E41").Value = '  -0.47%  '
$ws.Range("E42").Value = '  +0.03%  '
$ws.Range("D43").Value = "'25.16"
$ws.Range("E43").Value = '  -5.63%  '
$ws.Range("E44").Value = '  -0.05%  '
$ws.Range("E45").Value = '  +3.46%  '
$ws.Range("E46").Value = '  -3.16%  '
$ws.Range("E47").Value = '  -0.65%  '
$ws.Range("D48").Value = '2.376.38'
$ws.Range("E48").Value = '  -4.30%  '
$ws.Range("E49").Value = '  -2.50%  '
$ws.Range("E50").Value = '  -1.45%  '
$ws.Range("D51").Value = "'21.05"
